# Refresh tracking / rate / result data in rows 2-5 (DEV URL configuration change)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("P2") "320018787850"

Set-TextValue $ws.Range("P3") "320018787860"

Set-TextValue $ws.Range("P4") "320018764881"
Set-TextValue $ws.Range("Q4") '$49.70'
$ws.Range("R4").Value = "PASS"

Set-TextValue $ws.Range("P5") "320018766656"
